$wb = $excel.ActiveWorkbook

# --- Create the two new worksheets -----------------------------------------
# Existing workbook currently only has "AddCustomerTest".
# Target final sheet order (left to right): test_suite, AddCustomerTest, OpenAccountTest

$newOpen = $wb.Worksheets.Add()
$newOpen.Name = "OpenAccountTest"

$newSuite = $wb.Worksheets.Add()
$newSuite.Name = "test_suite"

# Put test_suite right before OpenAccountTest -> order: test_suite, OpenAccountTest, AddCustomerTest
$newSuite.Move($newOpen)

# Move AddCustomerTest to sit right after test_suite -> order: test_suite, AddCustomerTest, OpenAccountTest
$addCust = $wb.Worksheets.Item("AddCustomerTest")
$addCust.Move($null, $newSuite)

# --- Fill in "AddCustomerTest" sheet ----------------------------------------
$ws = $wb.Worksheets.Item("AddCustomerTest")

$ws.Cells.Item(1,1).Value = "firstname"
$ws.Cells.Item(1,2).Value = "lastname"
$ws.Cells.Item(1,3).Value = "postcode"
$ws.Cells.Item(1,4).Value = "alerttext"
$ws.Cells.Item(1,5).Value = "runmode"

$ws.Cells.Item(2,1).Value = "Raman"
$ws.Cells.Item(2,2).Value = "Arora"
$ws.Cells.Item(2,3).Value = 33518
$ws.Cells.Item(2,4).Value = "Customer added successfully"
$ws.Cells.Item(2,5).Value = "Y"

$ws.Cells.Item(3,1).Value = "Rahul"
$ws.Cells.Item(3,2).Value = "Bam"
$ws.Cells.Item(3,3).Value = 33518
$ws.Cells.Item(3,4).Value = "Customer added successfully"
$ws.Cells.Item(3,5).Value = "N"

$ws.Cells.Item(4,1).Value = "Andrew"
$ws.Cells.Item(4,2).Value = "Malkov"
$ws.Cells.Item(4,3).Value = 33518
$ws.Cells.Item(4,4).Value = "Customer added successfully"
$ws.Cells.Item(4,5).Value = "Y"

$ws.Cells.Item(5,1).Value = "Jack"
$ws.Cells.Item(5,2).Value = "Bim"
$ws.Cells.Item(5,3).Value = 33518
$ws.Cells.Item(5,4).Value = "Customer added successfully"
$ws.Cells.Item(5,5).Value = "Y"

$ws.Range("I10").Select() | Out-Null

# --- Fill in "test_suite" sheet ---------------------------------------------
$ws2 = $wb.Worksheets.Item("test_suite")

$ws2.Cells.Item(1,1).Value = "TCID"
$ws2.Cells.Item(1,2).Value = "runmode"

$ws2.Cells.Item(2,1).Value = "AddCustomerTest"
$ws2.Cells.Item(2,2).Value = "Y"

$ws2.Cells.Item(3,1).Value = "BankManagerLoginTest"
$ws2.Cells.Item(3,2).Value = "Y"

$ws2.Cells.Item(4,1).Value = "OpenAccountTest"
$ws2.Cells.Item(4,2).Value = "N"

$ws2.Columns.Item(1).AutoFit()

$ws2.Range("E13").Select() | Out-Null

# --- Fill in "OpenAccountTest" sheet ----------------------------------------
$ws3 = $wb.Worksheets.Item("OpenAccountTest")

$ws3.Cells.Item(1,1).Value = "customer"
$ws3.Cells.Item(1,2).Value = "currency"

$ws3.Cells.Item(2,1).Value = "Raman Arora"
$ws3.Cells.Item(2,2).Value = "Rupee"

$ws3.Columns.Item(1).AutoFit()

$ws3.Range("A3").Select() | Out-Null

# --- Make AddCustomerTest the active tab ------------------------------------
$wb.Worksheets.Item("AddCustomerTest").Activate()
